$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 946.8570999999999
$ws.Range("I40").Value = 940.5
$ws.Range("K40").Value = 940.5
$ws.Range("M40").Value = -765.5

$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3752
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -3142
$ws.Range("N67").ClearContents()

$ws.Range("H135").Value = 766.3488
$ws.Range("I135").Value = 546.5897
$ws.Range("J135").Value = 2909
$ws.Range("K135").Value = 4919.3073
$ws.Range("L135").Value = 26181
$ws.Range("M135").Value = -2384.3073
$ws.Range("N135").Value = -31251

$ws.Range("H137").Value = 2162.1167
$ws.Range("I137").Value = 2184.2827
$ws.Range("J137").Value = 2089.2856
$ws.Range("K137").Value = 6552.848100000001
$ws.Range("L137").Value = 6267.8568
$ws.Range("M137").Value = -4002.848100000001
$ws.Range("N137").Value = -11367.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 999.08
$ws.Range("I32").Value = 950.5979599999999
$ws.Range("K32").Value = 950.5979599999999
$ws.Range("M32").Value = -663.5979599999999

$ws.Range("H47").Value = 28500
$ws.Range("J47").Value = 28500
$ws.Range("L47").Value = 28500
$ws.Range("N47").Value = -29950

$ws.Range("H74").Value = 781.5122
$ws.Range("I74").Value = 684.89655
$ws.Range("J74").Value = 1015
$ws.Range("K74").Value = 684.89655
$ws.Range("L74").Value = 1015
$ws.Range("M74").Value = 189.10345
$ws.Range("N74").Value = -2763

$ws.Range("H77").Value = 781.5122
$ws.Range("I77").Value = 684.89655
$ws.Range("J77").Value = 1015
$ws.Range("K77").Value = 3424.48275
$ws.Range("L77").Value = 5075
$ws.Range("M77").Value = 943.5172499999999
$ws.Range("N77").Value = -13811

$ws.Range("H122").Value = 2065.5454
$ws.Range("I122").Value = 1607.3914
$ws.Range("J122").Value = 3119.3
$ws.Range("K122").Value = 4822.174199999999
$ws.Range("L122").Value = 9357.900000000001
$ws.Range("M122").Value = -2372.174199999999
$ws.Range("N122").Value = -14257.9

$ws.Range("H132").Value = 26318412
$ws.Range("I132").Value = 38463332
$ws.Range("J132").Value = 4417.6665
$ws.Range("K132").Value = 115389996
$ws.Range("L132").Value = 13252.9995
$ws.Range("M132").Value = -115387466
$ws.Range("N132").Value = -18312.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1765.5625
$ws.Range("I105").Value = 1495
$ws.Range("J105").Value = 2036.125
$ws.Range("K105").Value = 1495
$ws.Range("L105").Value = 2036.125
$ws.Range("M105").Value = 252
$ws.Range("N105").Value = -5530.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2905.75
$ws.Range("I16").Value = 1788.5
$ws.Range("J16").Value = 3278.1667
$ws.Range("K16").Value = 1788.5
$ws.Range("L16").Value = 3278.1667
$ws.Range("M16").Value = -1501.5
$ws.Range("N16").Value = -3852.1667

$ws.Range("H31").Value = 2864.6956
$ws.Range("I31").Value = 1856.2
$ws.Range("J31").Value = 3640.4614
$ws.Range("K31").Value = 1856.2
$ws.Range("L31").Value = 3640.4614
$ws.Range("M31").Value = -1561.2
$ws.Range("N31").Value = -4230.4614

$ws.Range("H34").Value = 2864.6956
$ws.Range("I34").Value = 1856.2
$ws.Range("J34").Value = 3640.4614
$ws.Range("K34").Value = 1856.2
$ws.Range("L34").Value = 3640.4614
$ws.Range("M34").Value = -1654.2
$ws.Range("N34").Value = -4044.4614

$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 30000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30588

$ws.Range("H62").Value = 4666.6665
$ws.Range("J62").Value = 4798.5713
$ws.Range("L62").Value = 4798.5713
$ws.Range("N62").Value = -6046.5713

$ws.Range("H65").Value = 4666.6665
$ws.Range("J65").Value = 4798.5713
$ws.Range("L65").Value = 23992.8565
$ws.Range("N65").Value = -30232.8565

$ws.Range("H113").Value = 2905.75
$ws.Range("I113").Value = 1788.5
$ws.Range("J113").Value = 3278.1667
$ws.Range("K113").Value = 1788.5
$ws.Range("L113").Value = 3278.1667
$ws.Range("M113").Value = 381.5
$ws.Range("N113").Value = -7618.1667

$ws.Range("H139").Value = 39850
$ws.Range("J139").Value = 39850
$ws.Range("L139").Value = 39850
$ws.Range("N139").Value = -50130

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 148.54167
$ws.Range("J12").Value = 225.13333
$ws.Range("L12").Value = 675.39999
$ws.Range("N12").Value = -1021.39999

$ws.Range("H107").Value = 1754.8334
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1754.8334
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 5264.5002
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -9104.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 62013.332
$ws.Range("J48").Value = 85030
$ws.Range("L48").Value = 85030
$ws.Range("N48").Value = -86000

$ws.Range("H102").Value = 49867.24
$ws.Range("I102").Value = 1307.0714
$ws.Range("J102").Value = 146987.58
$ws.Range("K102").Value = 1307.0714
$ws.Range("L102").Value = 146987.58
$ws.Range("M102").Value = 314.9286
$ws.Range("N102").Value = -150231.58

$ws.Range("H107").Value = 1043.2
$ws.Range("I107").Value = 497.91666
$ws.Range("J107").Value = 1546.5385
$ws.Range("K107").Value = 497.91666
$ws.Range("L107").Value = 1546.5385
$ws.Range("M107").Value = 1422.08334
$ws.Range("N107").Value = -5386.538500000001

$ws.Range("H132").Value = 2979.366
$ws.Range("I132").Value = 2726.8096
$ws.Range("J132").Value = 3244.55
$ws.Range("K132").Value = 8180.4288
$ws.Range("L132").Value = 9733.650000000001
$ws.Range("M132").Value = -5650.4288
$ws.Range("N132").Value = -14793.65

$ws.Range("H138").Value = 78000
$ws.Range("J138").Value = 78000
$ws.Range("L138").Value = 78000
$ws.Range("N138").Value = -88280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 924.58826
$ws.Range("I16").Value = 964.125
$ws.Range("J16").Value = 292
$ws.Range("K16").Value = 964.125
$ws.Range("L16").Value = 292
$ws.Range("M16").Value = -794.125
$ws.Range("N16").Value = -632

$ws.Range("H55").Value = 797.65
$ws.Range("I55").Value = 203.875
$ws.Range("J55").Value = 1193.5
$ws.Range("K55").Value = 203.875
$ws.Range("L55").Value = 1193.5
$ws.Range("M55").Value = -30.875
$ws.Range("N55").Value = -1539.5

$ws.Range("H122").Value = 2929.7026
$ws.Range("I122").Value = 2523.1785
$ws.Range("J122").Value = 4194.4443
$ws.Range("K122").Value = 7569.5355
$ws.Range("L122").Value = 12583.3329
$ws.Range("M122").Value = -5119.5355
$ws.Range("N122").Value = -17483.3329

$ws.Range("H132").Value = 2996.8235
$ws.Range("I132").Value = 1978.8
$ws.Range("K132").Value = 5936.4
$ws.Range("M132").Value = -3406.4

$ws.Range("H135").Value = 29671.727
$ws.Range("J135").Value = 29671.727
$ws.Range("L135").Value = 29671.727
$ws.Range("N135").Value = -39811.727

$ws.Range("H136").Value = 1384.8223
$ws.Range("I136").Value = 829.5806
$ws.Range("J136").Value = 2614.2856
$ws.Range("K136").Value = 2488.7418
$ws.Range("L136").Value = 7842.8568
$ws.Range("M136").Value = 61.25820000000022
$ws.Range("N136").Value = -12942.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 272255.28
$ws.Range("I122").Value = 346467.66
$ws.Range("J122").Value = 3235.375
$ws.Range("K122").Value = 1039402.98
$ws.Range("L122").Value = 9706.125
$ws.Range("M122").Value = -1036952.98
$ws.Range("N122").Value = -14606.125

$ws.Range("H127").Value = 35000
$ws.Range("J127").Value = 35000
$ws.Range("L127").Value = 35000
$ws.Range("N127").Value = -44920

$ws.Range("H132").Value = 8300.154
$ws.Range("I132").Value = 1767.875
$ws.Range("K132").Value = 5303.625
$ws.Range("M132").Value = -2773.625

$ws.Range("H136").Value = 946.9
$ws.Range("I136").Value = 577.3611
$ws.Range("J136").Value = 1897.1428
$ws.Range("K136").Value = 1732.0833
$ws.Range("L136").Value = 5691.428400000001
$ws.Range("M136").Value = 817.9167000000002
$ws.Range("N136").Value = -10791.4284
